# Add a new worksheet "Unfriend_Facebook_Friends" with instructions for
# how many friends to unfriend, and update the Friend_Request_Management
# sheet's Skip FRM controls.

$wb = $excel.ActiveWorkbook

# --- Update existing "Friend_Request_Management" sheet ---
$frm = $wb.Worksheets.Item("Friend_Request_Management")

# Clear old demo values in A2/B2 (no longer used)
[void]$frm.Range("A2").ClearContents()
[void]$frm.Range("B2").ClearContents()

# D3 gets the long explanatory note about the Skip FRM flag
$frm.Range("D3").Value = "Skip FRM: when you write yes the seq. will execute otherwise it will skipped."

# C2 now holds the "no" flag for Skip FRM
$frm.Range("C2").Value = "no"

# Move the selection as recorded after editing
[void]$frm.Range("B12").Select()

# --- Add new "Unfriend_Facebook_Friends" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Unfriend_Facebook_Friends"

$newSheet.Range("A1").Value = "No_Frineds_to_Remove"
$newSheet.Range("A2").Value = 5

$newSheet.Columns.Item(1).ColumnWidth = 22.7109375

[void]$newSheet.Range("A7").Select()
